# SignalAnalysis Greek (el-GR) translation workbook update
#
# Inserts a new translation row ("strFileHeader29") right before the
# existing "strFileHeaderSection" row (the row that used to be row 50),
# pushing every following row down by one. The new row's Comment/English
# cells use a left/vertical-center alignment (no wrap), and the Comment +
# English cells are additionally unlocked (matches the authored xf with
# applyProtection + protection locked="0").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The translation table lives on the only worksheet, backed by a single
# Excel Table ("Tabla13") that currently spans B2:E158.
$lo = $ws.ListObjects.Item(1)

# Physically insert a blank row at sheet row 50 - this shifts all data
# (and row heights) in rows 50..158 down to 51..159, exactly like the
# diff shows.
$ws.Rows.Item(50).Insert()

# The table's own `ref`/`autoFilter` range doesn't auto-grow from a plain
# row insert, so resize it explicitly to include the new row.
$lo.Resize($ws.Range("B2:E159"))

# Populate the new row's three translation columns. New literal strings
# are appended to the shared-string table in the order they're written,
# which reproduces the new uniqueCount=310 / count=334 entries
# (strFileHeader29, "Field description in exported file",
# "Differentiation algorithm") at indices 307-309.
$ws.Range("B50").Value = "strFileHeader29"
$ws.Range("C50").Value = "Field description in exported file"
$ws.Range("D50").Value = "Differentiation algorithm"

# New row formatting: left/center alignment without wrap text (unlike the
# existing style used elsewhere in the sheet, which wraps).
$row50 = $ws.Range("B50:D50")
$row50.HorizontalAlignment = -4131 # xlLeft
$row50.VerticalAlignment = -4108   # xlCenter
$row50.WrapText = $false

# Comment (C50) and English (D50) cells are additionally unlocked.
$ws.Range("C50:D50").Locked = $false
